$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire _uuid column (column I), shifting everything right of it
# one column to the left.
$ws.Columns("I").Delete()

# Update the selection to match the post-edit state (cursor ends up on I4
# after the column delete / re-selection in the original edit).
$ws.Range("I4").Select()
